$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Database password (C6) was rotated to a new value
$ws.Range("C6").Value = "hgtunb26364AQK"

# Insert a new "注意点" (Notes) row below the URL row
$ws.Range("B8").Value = "注意点"
$ws.Range("C8").Value = "バーセルからのdb接続urlはtransaction poolerにする。DirctはNG"

# Give B8 the same label style used by the other field names in column B (B4:B7)
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Apply character-level rich-text formatting to the "注意点" label
# (applied as two adjoining runs so it is stored as shared-string run formatting)
$f1 = $ws.Range("B8").Characters(1, 2).Font
$f1.Name = "游ゴシック"
$f1.Color = 1513239
$f1.Size = 16
$f2 = $ws.Range("B8").Characters(3, 1).Font
$f2.Name = "游ゴシック"
$f2.Color = 1513239
$f2.Size = 16

# Row grew taller to fit the bigger "注意点" label font
$ws.Range("B8").RowHeight = 25.5

# Leave the selection where the user was last working
$ws.Range("C9").Select()
